$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs ("<w:p/>").
# Both need to become justified, "Kohinoor Bangla Light" paragraphs;
# the last one additionally receives the "resposta 3" text.
$count = $d.Paragraphs.Count
$pBlank = $d.Paragraphs.Item($count - 1)
$pText  = $d.Paragraphs.Item($count)

$font = "Kohinoor Bangla Light"

# --- blank paragraph: justify + set the (empty) paragraph-mark font ---
$pBlank.Alignment = 3
$pBlank.Range.Font.NameAscii = $font
$pBlank.Range.Font.NameOther = $font
$pBlank.Range.Font.NameBi = $font

# --- last paragraph: justify + paragraph-mark font, then the new text ---
$pText.Alignment = 3
$pText.Range.Font.NameAscii = $font
$pText.Range.Font.NameOther = $font
$pText.Range.Font.NameBi = $font

$texto = "3. Como comentado anteriormente" + `
    " de acordo com as necessidades do cliente" + `
    " o método ágil se encaixa muito bem neste projeto. Ele garante a entrega de pequenas partes do produto final da maneira mais rápida possível mantendo a qualidade desejada" + `
    ". Possibilitando a compreensão das “dores” do c" + `
    "liente, a comunicação constante com as partes envolvidas" + `
    ", entregas no final de cada Sprint" + `
    " e o trabalho independente entre as equipes" + `
    "."

$pText.Range.InsertAfter($texto)

# Re-apply the font to the run that was just inserted.
$pText.Range.Font.NameAscii = $font
$pText.Range.Font.NameOther = $font
$pText.Range.Font.NameBi = $font
$pText.Alignment = 3
